$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "project"
$ws.Range("B1").Value = "BMI"
$ws.Range("C1").Value = "change proneness"

# --- Column A: renamed project/version labels ---
$ws.Range("A2").Value  = "Apache commons Lang 3.0-3.6"
$ws.Range("A3").Value  = "Apache commons Lang 3.6-3.7"
$ws.Range("A4").Value  = "Apache commons Lang 3.7-3.8"
$ws.Range("A5").Value  = "Apache commons codec 1.10-1.11"
$ws.Range("A6").Value  = "Apache commons codec 1.11-1.12"
$ws.Range("A7").Value  = "Apache commons codec 1.9-1.10"
$ws.Range("A8").Value  = "Apache commons collections 3.2-4.0"
$ws.Range("A9").Value  = "Apache commons collections 4.0-4.1"
$ws.Range("A10").Value = "Apache commons collections 4.1-4.3"
$ws.Range("A11").Value = "Apache commons configuration 2.1-2.2"
$ws.Range("A12").Value = "Apache commons configuration 2.2-2.3"
$ws.Range("A13").Value = "Apache commons configuration 2.3-2.4"
$ws.Range("A14").Value = "Jfreechart 1.0.18-1.0.19"
$ws.Range("A15").Value = "Jfreechart 0.19-1.5.0"

# --- Column B: BMI values rescaled from fraction to percent (x100) ---
# (multiplied from the original fractional values to keep identical
#  floating-point rounding artifacts to the source edit)
$ws.Range("B2").Value  = 0.10833 * 100
$ws.Range("B3").Value  = 0.43332999999999999 * 100
$ws.Range("B4").Value  = 0.33968300000000001 * 100
$ws.Range("B5").Value  = 0.30555599999999999 * 100
$ws.Range("B6").Value  = 0.44444400000000001 * 100
$ws.Range("B7").Value  = 1 * 100
$ws.Range("B8").Value  = 0.40317500000000001 * 100
$ws.Range("B9").Value  = 0.38611000000000001 * 100
$ws.Range("B10").Value = 0.41666700000000001 * 100
$ws.Range("B11").Value = 0.66666999999999998 * 100
$ws.Range("B12").Value = 0.15757599999999999 * 100
$ws.Range("B13").Value = 0.030303 * 100
$ws.Range("B14").Value = 2.5 * 100
$ws.Range("B15").Value = 0.66666999999999998 * 100

# --- Styling: column B uses the "applied number format" style ---
$ws.Range("B1:B15").NumberFormat = "General"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 43.5
$ws.Columns.Item(2).ColumnWidth = 10.83203125
$ws.Columns.Item(4).ColumnWidth = 10.83203125

# --- Selection cursor moved ---
$ws.Range("E6").Select()
